$d = $word.ActiveDocument

# Locate the paragraph containing the "Please read Using Pressure Canners..."
# text. We also remove the blank paragraph immediately preceding it (the
# empty line that separated it from the prior paragraph), so that the
# remaining blank-line spacing between paragraphs stays consistent.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Please read Using Pressure Canners*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    $target = $d.Paragraphs.Item($targetIndex)
    $blank = $d.Paragraphs.Item($targetIndex - 1)

    $startPos = $blank.Range.Start
    $endPos = $target.Range.End

    $r = $d.Range($startPos, $endPos)
    $r.Delete()
}
